# Updates the crypto price/volume snapshot table (columns B:E, rows 2-51)
# to the new scrape values. Values in column D (Price) are written as text
# (via NumberFormat "@") whenever they could otherwise be auto-coerced into
# numbers by Excel (e.g. "1.000" -> 1), to preserve the original formatting
# exactly as scraped. Rows 41/42 additionally swap the Coin name and Link
# (RenderToken <-> MXToken) to match the refreshed ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.384.17"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.828.19"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.09"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4451"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3769"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07409"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8788"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.88"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "1.829.75"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.721"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.436"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.08"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07060"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008813"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.10"
$ws.Range("D21").Value = "27.392.70"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.358"
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.952"
$ws.Range("E24").Value = "  -1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.27"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.351"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.10"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08906"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7965"
$ws.Range("E31").Value = "  +6.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.200"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.560"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.966"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9997"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.106"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01980"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05275"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.298"
$ws.Range("E39").Value = "  +3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5338"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.875"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.342"
$ws.Range("E42").Value = "  +18.53%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.713"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5074"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.66"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.50"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.688"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9996"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06394"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.22"
$ws.Range("E51").Value = "  +5.75%  "
